# Bulk refresh of market-price derived columns (H:N) across all job sheets.
# Source data regenerated by the scheduled pricing runner; values below are the
# newly computed figures for each affected leve row (keyed by its row number).
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 138.23077
$ws.Range("I9").Value = 138.45454
$ws.Range("K9").Value = 138.45454
$ws.Range("M9").Value = 30.54545999999999
# Row 32
$ws.Range("H32").Value = 1703.3334
$ws.Range("I32").Value = 1158.3334
$ws.Range("J32").Value = 2066.6667
$ws.Range("K32").Value = 1158.3334
$ws.Range("L32").Value = 2066.6667
$ws.Range("M32").Value = -832.3334
$ws.Range("N32").Value = -2718.6667
# Row 62
$ws.Range("H62").Value = 2516.6667
$ws.Range("I62").Value = 2516.6667
$ws.Range("K62").Value = 2516.6667
$ws.Range("M62").Value = -1892.6667
# Row 65
$ws.Range("H65").Value = 2516.6667
$ws.Range("I65").Value = 2516.6667
$ws.Range("K65").Value = 12583.3335
$ws.Range("M65").Value = -9463.333500000001
# Row 80
$ws.Range("H80").Value = 878.82355
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 621.25
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 1863.75
$ws.Range("M80").Value = -14002
$ws.Range("N80").Value = -3859.75
# Row 83
$ws.Range("H83").Value = 878.82355
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 621.25
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 5591.25
$ws.Range("M83").Value = -40008
$ws.Range("N83").Value = -15575.25
# Row 111
$ws.Range("H111").Value = 3613.2856
$ws.Range("I111").Value = 3847
$ws.Range("K111").Value = 11541
$ws.Range("M111").Value = -8474
# Row 118
$ws.Range("H118").Value = 400
$ws.Range("I118").Value = 400
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1200
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 457
$ws.Range("N118").Value = $null
# Row 132
$ws.Range("H132").Value = 7757962.5
$ws.Range("I132").Value = 9806617
$ws.Range("J132").Value = 18600.666
$ws.Range("K132").Value = 29419851
$ws.Range("L132").Value = 55801.99800000001
$ws.Range("M132").Value = -29417321
$ws.Range("N132").Value = -60861.99800000001
# Row 137
$ws.Range("H137").Value = 1140.826
$ws.Range("I137").Value = 1160.9412
$ws.Range("J137").Value = 1083.8334
$ws.Range("K137").Value = 3482.8236
$ws.Range("L137").Value = 3251.5002
$ws.Range("M137").Value = -932.8235999999997
$ws.Range("N137").Value = -8351.5002

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4343.6855
$ws.Range("I32").Value = 4829.967
$ws.Range("K32").Value = 4829.967
$ws.Range("M32").Value = -4542.967
# Row 45
$ws.Range("H45").Value = 1687.5
$ws.Range("I45").Value = 1876.5
$ws.Range("J45").Value = 931.5
$ws.Range("K45").Value = 1876.5
$ws.Range("L45").Value = 931.5
$ws.Range("M45").Value = -1499.5
$ws.Range("N45").Value = -1685.5
# Row 74
$ws.Range("H74").Value = 1166.4375
$ws.Range("I74").Value = 838.63635
$ws.Range("K74").Value = 838.63635
$ws.Range("M74").Value = 35.36365000000001
# Row 77
$ws.Range("H77").Value = 1166.4375
$ws.Range("I77").Value = 838.63635
$ws.Range("K77").Value = 4193.18175
$ws.Range("M77").Value = 174.8182500000003
# Row 97
$ws.Range("H97").Value = 547.1429000000001
$ws.Range("I97").Value = 388.33334
$ws.Range("K97").Value = 388.33334
$ws.Range("M97").Value = 107.66666
# Row 125
$ws.Range("H125").Value = 44799.5
$ws.Range("J125").Value = 44799.5
$ws.Range("L125").Value = 44799.5
$ws.Range("N125").Value = -54639.5

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 539.13336
$ws.Range("I80").Value = 398
$ws.Range("J80").Value = 574.4167
$ws.Range("K80").Value = 398
$ws.Range("L80").Value = 574.4167
$ws.Range("M80").Value = 600
$ws.Range("N80").Value = -2570.4167
# Row 83
$ws.Range("H83").Value = 539.13336
$ws.Range("I83").Value = 398
$ws.Range("J83").Value = 574.4167
$ws.Range("K83").Value = 1990
$ws.Range("L83").Value = 2872.0835
$ws.Range("M83").Value = 3002
$ws.Range("N83").Value = -12856.0835
# Row 94
$ws.Range("H94").Value = 16667441
$ws.Range("I94").Value = 17857830
$ws.Range("K94").Value = 17857830
$ws.Range("M94").Value = -17857379
# Row 99
$ws.Range("H99").Value = 31251314
$ws.Range("I99").Value = 38462570
$ws.Range("J99").Value = 2533.3333
$ws.Range("K99").Value = 38462570
$ws.Range("L99").Value = 2533.3333
$ws.Range("M99").Value = -38461072
$ws.Range("N99").Value = -5529.3333

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2353.1333
$ws.Range("I31").Value = 1124.625
$ws.Range("K31").Value = 1124.625
$ws.Range("M31").Value = -829.625
# Row 34
$ws.Range("H34").Value = 2353.1333
$ws.Range("I34").Value = 1124.625
$ws.Range("K34").Value = 1124.625
$ws.Range("M34").Value = -922.625
# Row 94
$ws.Range("H94").Value = 1637.1666
# Row 132
$ws.Range("H132").Value = 6421.909
$ws.Range("I132").Value = 8729.538
$ws.Range("J132").Value = 3088.6667
$ws.Range("K132").Value = 26188.614
$ws.Range("L132").Value = 9266.000100000001
$ws.Range("M132").Value = -23658.614
$ws.Range("N132").Value = -14326.0001

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1099994.8
$ws.Range("I4").Value = 99993.11
$ws.Range("K4").Value = 299979.33
$ws.Range("M4").Value = -299867.33
# Row 5
$ws.Range("H5").Value = 550.3182
$ws.Range("I5").Value = 547.85
$ws.Range("K5").Value = 1643.55
$ws.Range("M5").Value = -1531.55
# Row 23
$ws.Range("H23").Value = 591.3333
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 309.6
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 928.8000000000001
$ws.Range("M23").Value = -5765
$ws.Range("N23").Value = -1398.8
# Row 32
$ws.Range("H32").Value = 2120
$ws.Range("J32").Value = 2120
$ws.Range("L32").Value = 6360
$ws.Range("N32").Value = -6926
# Row 113
$ws.Range("H113").Value = 755.375
$ws.Range("I113").Value = 760
$ws.Range("J113").Value = 753.2727
$ws.Range("K113").Value = 2280
$ws.Range("L113").Value = 2259.8181
$ws.Range("M113").Value = -110
$ws.Range("N113").Value = -6599.8181
# Row 122
$ws.Range("H122").Value = 864.7778
$ws.Range("I122").Value = 805.5
$ws.Range("K122").Value = 7249.5
$ws.Range("M122").Value = -4799.5
# Row 135
$ws.Range("H135").Value = 550.3182
$ws.Range("I135").Value = 547.85
$ws.Range("K135").Value = 4930.650000000001
$ws.Range("M135").Value = -2395.650000000001

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 7015000
$ws.Range("J11").Value = 8000000
$ws.Range("L11").Value = 8000000
$ws.Range("N11").Value = -8000278
# Row 49
$ws.Range("H49").Value = 18500
$ws.Range("J49").Value = 18500
$ws.Range("L49").Value = 18500
$ws.Range("N49").Value = -18868

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1981.3636
$ws.Range("I7").Value = 2051.6
$ws.Range("J7").Value = 1922.8334
$ws.Range("K7").Value = 2051.6
$ws.Range("L7").Value = 1922.8334
$ws.Range("M7").Value = -1939.6
$ws.Range("N7").Value = -2146.8334
# Row 22
$ws.Range("H22").Value = 1307.25
$ws.Range("J22").Value = 1671.6
$ws.Range("L22").Value = 1671.6
$ws.Range("N22").Value = -2261.6
# Row 27
$ws.Range("H27").Value = 1307.25
$ws.Range("J27").Value = 1671.6
$ws.Range("L27").Value = 1671.6
$ws.Range("N27").Value = -1885.6
# Row 42
$ws.Range("H42").Value = 18000
$ws.Range("J42").Value = 18000
$ws.Range("L42").Value = 18000
$ws.Range("N42").Value = -19126
# Row 49
$ws.Range("H49").Value = 18000
$ws.Range("J49").Value = 18000
$ws.Range("L49").Value = 18000
$ws.Range("N49").Value = -18294
# Row 69
$ws.Range("H69").Value = 67142.86
$ws.Range("J69").Value = 67142.86
$ws.Range("L69").Value = 67142.86
$ws.Range("N69").Value = -68764.86
# Row 72
$ws.Range("H72").Value = 67142.86
$ws.Range("J72").Value = 67142.86
$ws.Range("L72").Value = 201428.58
$ws.Range("N72").Value = -209540.58
# Row 100
$ws.Range("H100").Value = 1851.6875
$ws.Range("I100").Value = 1617.4615
$ws.Range("J100").Value = 2866.6667
$ws.Range("K100").Value = 1617.4615
$ws.Range("L100").Value = 2866.6667
$ws.Range("M100").Value = -1076.4615
$ws.Range("N100").Value = -3948.6667
# Row 126
$ws.Range("H126").Value = 1981.3636
$ws.Range("I126").Value = 2051.6
$ws.Range("J126").Value = 1922.8334
$ws.Range("K126").Value = 6154.799999999999
$ws.Range("L126").Value = 5768.5002
$ws.Range("M126").Value = -3684.799999999999
$ws.Range("N126").Value = -10708.5002

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 3000
# Row 4
$ws.Range("H4").Value = 2120
$ws.Range("J4").Value = 2120
$ws.Range("L4").Value = 2120
$ws.Range("N4").Value = -2346
# Row 8
$ws.Range("H8").Value = 15000
$ws.Range("I8").Value = 15000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -14860
$ws.Range("N8").Value = $null
# Row 19
$ws.Range("H19").Value = 899.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 899.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 899.5
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = -1247.5
# Row 21
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10470
# Row 29
$ws.Range("H29").Value = 1950
$ws.Range("J29").Value = 1950
$ws.Range("L29").Value = 1950
$ws.Range("N29").Value = -2530
# Row 35
$ws.Range("H35").Value = 10000
$ws.Range("J35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("N35").Value = -10580
# Row 86
$ws.Range("H86").Value = 14800
$ws.Range("J86").Value = 14800
$ws.Range("L86").Value = 14800
$ws.Range("N86").Value = -17046
# Row 89
$ws.Range("H89").Value = 14800
$ws.Range("J89").Value = 14800
$ws.Range("L89").Value = 74000
$ws.Range("N89").Value = -85232
# Row 126
$ws.Range("H126").Value = 52910804
$ws.Range("I126").Value = 65360030
$ws.Range("J126").Value = 1582.75
$ws.Range("K126").Value = 196080090
$ws.Range("L126").Value = 4748.25
$ws.Range("M126").Value = -196077620
$ws.Range("N126").Value = -9688.25
